$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Helper: isolate the text covered by $rng into its own <w:r> run
# without altering its visible formatting. Toggling Bold on then back
# off forces a run break at both edges of $rng; since the property is
# restored to its original value, the resulting runs end up with
# identical <w:rPr> to their neighbours but stay distinct <w:r>
# elements - i.e. exactly the segmentation Word leaves behind after it
# wraps a word with <w:proofErr w:type="spellStart"/> ...
# <w:proofErr w:type="spellEnd"/> during an interactive spell-check.
# ---------------------------------------------------------------------
function Split-Run($rng) {
    $rng.Font.Bold = $true
    $rng.Font.Bold = $false
}

# Find the next occurrence of $text at/after absolute offset $from,
# isolate it into its own run, and return the absolute end offset of
# the match (so callers can keep searching forward in document
# order).
function Isolate-NextMatch($d, $text, $from) {
    $rng = $d.Range($from, $d.Content.End)
    $ok = $rng.Find.Execute($text, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        throw "Could not find '$text' starting at $from"
    }
    Split-Run $rng
    return $rng.End
}

# ===== 1. "...but the localstorage difficulty will change things." =====
$pos = 0
$pos = Isolate-NextMatch $d "localstorage" $pos

# ===== 2. Casual CONCEPT TEXT paragraph: "...named Dealle, found out
#          ... half yokai, in some magical way ... the world. Dealle
#          was kept..." ================================================
$pos = Isolate-NextMatch $d "Dealle" $pos
$pos = Isolate-NextMatch $d "yokai" $pos

# Split "...destroys the world." | " " | "Dealle" | " was kept..." so
# the space and the name each become their own run.
$anchor = $d.Range($pos, $d.Content.End)
$anchor.Find.Execute("the world. Dealle", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$dotEnd = $anchor.Start + ("the world.").Length
$spaceEnd = $dotEnd + 1
$dealleEnd = $spaceEnd + ("Dealle").Length
Split-Run ($d.Range($dotEnd, $spaceEnd))
Split-Run ($d.Range($spaceEnd, $dealleEnd))
$pos = $dealleEnd

# ===== 3. IN-GAME-TEXT paragraph 1: "...witch named Dealle discovered
#          him... As Dealle tended..." ==================================
$pos = Isolate-NextMatch $d "Dealle" $pos
$pos = Isolate-NextMatch $d "Dealle" $pos

# ===== 4. IN-GAME-TEXT paragraph 2: "...the essence of a yokai,
#          through mystical means..." ===================================
$pos = Isolate-NextMatch $d "yokai" $pos

# ===== 5. IN-GAME-TEXT paragraph 4: "...perilous predicament, Dealle
#          found herself torn..." =======================================
$pos = Isolate-NextMatch $d "Dealle" $pos

# ===== 6. New blank paragraph right after "One random card is sent to
#          the used pile." and before the "Mechanics" heading. =========
$rng = $d.Content.Duplicate
$rng.Find.Execute("One random card is sent to the used pile.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Collapse(0)
$rng.InsertParagraphAfter()

Write-Host "done"
